$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The medication list in column A (rows 2:144) is kept sorted. A new
# medication "아모딘정" needs to be inserted (alphabetically) at row 80,
# so every existing entry from row 80 through row 143 shifts down by
# one row (row 144, previously blank, receives the old row 143 value).
for ($r = 144; $r -ge 81; $r--) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r - 1, 1).Value2
}
$ws.Cells.Item(80, 1).Value = "아모딘정"
